# LOM3097.xlsx update
# - Fills in the "Objetivos" text, adds a second responsible professor row,
#   and fills in several previously-empty syllabus/criteria/bibliography cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 10 ("Objetivos:") — replace the placeholder professor name with the
#    real objectives text.
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "Complementar a formação do estudante de Engenharia de Materiais proporcionando a oportunidade de desenvolver estágio não obrigatório para aprofundamento da experiência e vivência profissional em ambiente industrial"
$ws.Range("C10").Value = "Complementar a formação do estudante de Engenharia de Materiais proporcionando a oportunidade de desenvolver estágio não obrigatório para aprofundamento da experiência e vivência profissional em ambiente industrial"

# ---------------------------------------------------------------------------
# 2. Insert two new rows right after row 12 ("Docentes responsáveis:") to
#    hold the two professors, one per row, in columns B/C only.
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

# The inserted rows copy row 12's formatting (bold, col-A style) by default;
# column A must stay empty on these rows, so drop that cell entirely.
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()

# Pull correct B/C formatting (wrap-text body style / red "changed" style)
# from an existing, fully-populated data row before writing values in.
$ws.Range("B9").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C13").PasteSpecial(-4122)

$ws.Range("B9").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C14").PasteSpecial(-4122)

$ws.Range("B13").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("C13").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("B14").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C14").Value = "519033 - Carlos Yujiro Shigue"

# ---------------------------------------------------------------------------
# 3. Old row 13 ("Programa resumido:") is now row 15 after the insert — give
#    it its real short-syllabus text (row height is already 60 from the
#    original row and needs no change).
# ---------------------------------------------------------------------------
$ws.Range("B15").Value = "Elaboração do plano de trabalho de estágio. Realização do estágio. Elaboração de relatório final de estágio."
$ws.Range("C15").Value = "Elaboração do plano de trabalho de estágio. Realização do estágio. Elaboração de relatório final de estágio."

# ---------------------------------------------------------------------------
# 4. Old row 15 ("Programa:") is now row 17 — give it the real syllabus text.
# ---------------------------------------------------------------------------
$ws.Range("B17").Value = "O estágio será realizado sob a supervisão de docente designado pelo Coordenador de Estágio do curso de Engenharia de Materiais. O conteúdo será estabelecido individualmente no Plano de Trabalho entre o supervisor responsável pelo Estágio e o docente supervisor, desde que relacionado com as áreas afins da Engenharia de Materiais em ambiente de trabalho industrial. Apresentação de relatório sobre as atividades desenvolvidas no estágio."
$ws.Range("C17").Value = "O estágio será realizado sob a supervisão de docente designado pelo Coordenador de Estágio do curso de Engenharia de Materiais. O conteúdo será estabelecido individualmente no Plano de Trabalho entre o supervisor responsável pelo Estágio e o docente supervisor, desde que relacionado com as áreas afins da Engenharia de Materiais em ambiente de trabalho industrial. Apresentação de relatório sobre as atividades desenvolvidas no estágio."

# ---------------------------------------------------------------------------
# 5. Old row 18 ("Método:") is now row 20 — give it the real method text.
# ---------------------------------------------------------------------------
$ws.Range("B20").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."
$ws.Range("C20").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."

# ---------------------------------------------------------------------------
# 6. Old row 19 ("Critério:") is now row 21 — give it the real criteria text.
# ---------------------------------------------------------------------------
$ws.Range("B21").Value = "Critério`nMF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo docente supervisor do estágio."
$ws.Range("C21").Value = "Critério`nMF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo docente supervisor do estágio."

# ---------------------------------------------------------------------------
# 7. Old row 20 ("Norma de recuperação:") is now row 22 — give it its text.
# ---------------------------------------------------------------------------
$ws.Range("B22").Value = "Não será oferecida recuperação."
$ws.Range("C22").Value = "Não será oferecida recuperação."

# ---------------------------------------------------------------------------
# 8. Old row 21 ("Bibliografia:") is now row 23 — give it its text.
# ---------------------------------------------------------------------------
$ws.Range("B23").Value = "A ser definida com o supervisor responsável pelo estágio e pelo docente orientador em função das atividades desenvolvidas no estágio."
$ws.Range("C23").Value = "A ser definida com o supervisor responsável pelo estágio e pelo docente orientador em função das atividades desenvolvidas no estágio."
